# Update course promotion sheet:
#  - course/department name changed from "EDISON SCHOOL OF TECH SCIENCES" to "Automotive"
#  - row 2 height reduced to fit the shorter text
#  - selection left on the edited cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Automotive"
$ws.Rows.Item(2).RowHeight = 57
[void]$ws.Range("C2").Select()
